$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new data rows are inserted into the table right before the current
# row 572 (which holds the 2021-10-25 "1a nueva(o)" record). This pushes the
# existing rows 572-646 down to 574-648, matching the new dimension
# A1:R648.
$ws.Rows.Item(572).Insert()
$ws.Rows.Item(572).Insert()

# Fill in the fixed (repeated) columns for the two new rows, matching the
# rest of the table (same market / region / category / classification).
$ws.Cells.Item(572, 1).Value = 5
$ws.Cells.Item(572, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(572, 3).Value = "Maule"
$ws.Cells.Item(572, 5).Value = 7
$ws.Cells.Item(572, 6).Value = 100112004
$ws.Cells.Item(572, 7).Value = "Cebolla"
$ws.Cells.Item(572, 8).Value = "Sin especificar"
$ws.Cells.Item(572, 18).Value = "Hortaliza"

$ws.Cells.Item(573, 1).Value = 5
$ws.Cells.Item(573, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(573, 3).Value = "Maule"
$ws.Cells.Item(573, 5).Value = 7
$ws.Cells.Item(573, 6).Value = 100112004
$ws.Cells.Item(573, 7).Value = "Cebolla"
$ws.Cells.Item(573, 8).Value = "Sin especificar"
$ws.Cells.Item(573, 18).Value = "Hortaliza"

# New row 572: $/malla 15 kilos record for 2022-08-03
$ws.Cells.Item(572, 4).Value = 44776
$ws.Cells.Item(572, 9).Value = "1a (guarda)"
$ws.Cells.Item(572, 10).Value = 2000
$ws.Cells.Item(572, 11).Value = 5500
$ws.Cells.Item(572, 12).Value = 5500
$ws.Cells.Item(572, 13).Value = 5500
$ws.Cells.Item(572, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(572, 15).Value = "Región del Maule"
$ws.Cells.Item(572, 16).Value = 367
$ws.Cells.Item(572, 17).Value = 15

# New row 573: $/malla 25 kilos record for 2022-08-03
$ws.Cells.Item(573, 4).Value = 44776
$ws.Cells.Item(573, 9).Value = "1a (guarda)"
$ws.Cells.Item(573, 10).Value = 2500
$ws.Cells.Item(573, 11).Value = 8000
$ws.Cells.Item(573, 12).Value = 8000
$ws.Cells.Item(573, 13).Value = 8000
$ws.Cells.Item(573, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(573, 15).Value = "Región del Maule"
$ws.Cells.Item(573, 16).Value = 320
$ws.Cells.Item(573, 17).Value = 25
